$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.789999999999999
$ws.Range("D6").Value = -8.173999999999999
$ws.Range("D7").Value = -7.585000000000001
$ws.Range("E7").Value = 16.803
$ws.Range("D8").Value = -8.301
$ws.Range("E11").Value = 17.271
$ws.Range("E12").Value = 17.694
$ws.Range("E15").Value = 15.955
$ws.Range("D16").Value = -8.625999999999999
$ws.Range("D20").Value = -7.797
$ws.Range("E20").Value = 16.457
$ws.Range("D21").Value = -8.322000000000001
$ws.Range("E21").Value = 16.373
$ws.Range("E22").Value = 16.506
$ws.Range("E23").Value = 16.434
$ws.Range("D28").Value = -7.887
$ws.Range("D29").Value = -7.56
$ws.Range("E29").Value = 16.64
$ws.Range("D30").Value = -7.185
$ws.Range("D32").Value = -7.665000000000001
$ws.Range("E34").Value = 16.729
$ws.Range("D40").Value = -8.119000000000002
$ws.Range("E42").Value = 16.539
$ws.Range("E43").Value = 16.896
$ws.Range("E44").Value = 16.518
$ws.Range("E45").Value = 16.801
$ws.Range("D46").Value = -8.022
$ws.Range("E46").Value = 16.906
$ws.Range("E50").Value = 16.265
$ws.Range("D51").Value = -8.241000000000001
$ws.Range("E51").Value = 16.489
$ws.Range("D52").Value = -7.87
$ws.Range("D57").Value = -7.946000000000001
$ws.Range("E57").Value = 16.582
$ws.Range("D59").Value = -8.124000000000001
$ws.Range("D62").Value = -7.946
$ws.Range("E65").Value = 17.109
$ws.Range("D66").Value = -7.102000000000001
$ws.Range("E66").Value = 17.325
$ws.Range("E67").Value = 17.227
$ws.Range("D73").Value = -7.903
$ws.Range("D74").Value = -7.883999999999999
$ws.Range("D77").Value = -7.840000000000001
$ws.Range("E79").Value = 16.957
$ws.Range("E84").Value = 16.606
$ws.Range("E87").Value = 16.603
$ws.Range("D92").Value = -7.401999999999999
$ws.Range("E92").Value = 16.832
$ws.Range("E97").Value = 16.787
$ws.Range("D100").Value = -8.295
